# Add playlist and song detail
# - renumber the continuous "no" column on the "songs" sheet (rows 8-40)
# - update the remembered selection/active-cell on both sheets

$wb = $excel.ActiveWorkbook

$albums = $wb.Worksheets.Item("albums")
$songs  = $wb.Worksheets.Item("songs")

# --- songs sheet: renumber column A (rows 8..40) to a continuous sequence ---
# row 7 already holds 6, so row 8 continues the count with 7, 8, 9, ... 39
for ($r = 8; $r -le 40; $r++) {
    $songs.Cells.Item($r, 1).Value = $r - 1
}

# --- view state: selection / active cell on each sheet ---
$albums.Activate()
[void]$albums.Range("B8").Select()

$songs.Activate()
[void]$songs.Range("B3").Select()
